$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement trial_total (F) by 106 for rows 2-42, and update distractor (L:V) blocks
# to rebalance n distractor = n targets (per commit message).

$ws.Range("F2").Value = 122
$ws.Range("F3").Value = 123
$ws.Range("F4").Value = 124
$ws.Range("L4").Value = 'stimuli/img_ua9bs.png'
$ws.Range("M4").Value = 82
$ws.Range("N4").Value = 62.23333333333333
$ws.Range("O4").Value = 72.11666666666667
$ws.Range("P4").Value = 30
$ws.Range("Q4").Value = 9
$ws.Range("R4").Value = 9
$ws.Range("S4").Value = 9
$ws.Range("T4").Value = 9
$ws.Range("U4").Value = 9
$ws.Range("V4").Value = 9
$ws.Range("F5").Value = 125
$ws.Range("F6").Value = 126
$ws.Range("F7").Value = 127
$ws.Range("F8").Value = 128
$ws.Range("F9").Value = 129
$ws.Range("F10").Value = 130
$ws.Range("F11").Value = 131
$ws.Range("L11").Value = 'stimuli/img_a8wvq.png'
$ws.Range("M11").Value = 86.25925925925925
$ws.Range("N11").Value = 66.25925925925925
$ws.Range("O11").Value = 76.25925925925925
$ws.Range("P11").Value = 27
$ws.Range("Q11").Value = 10
$ws.Range("R11").Value = 10
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 10
$ws.Range("U11").Value = 10
$ws.Range("V11").Value = 10
$ws.Range("F12").Value = 132
$ws.Range("L12").Value = 'stimuli/img_c0me7.png'
$ws.Range("M12").Value = 68.4
$ws.Range("N12").Value = 45.62857142857143
$ws.Range("O12").Value = 57.01428571428572
$ws.Range("P12").Value = 35
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 4
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 4
$ws.Range("U12").Value = 4
$ws.Range("V12").Value = 4
$ws.Range("F13").Value = 133
$ws.Range("L13").Value = 'stimuli/img_ifebc.png'
$ws.Range("M13").Value = 84
$ws.Range("N13").Value = 65.88235294117646
$ws.Range("O13").Value = 74.94117647058823
$ws.Range("P13").Value = 34
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = 9
$ws.Range("S13").Value = 9
$ws.Range("T13").Value = 9
$ws.Range("U13").Value = 9
$ws.Range("V13").Value = 9
$ws.Range("F14").Value = 134
$ws.Range("L14").Value = 'stimuli/img_xesl0.png'
$ws.Range("M14").Value = 69.28571428571429
$ws.Range("N14").Value = 47.35714285714285
$ws.Range("O14").Value = 58.32142857142857
$ws.Range("P14").Value = 28
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 5
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 5
$ws.Range("U14").Value = 5
$ws.Range("V14").Value = 5
$ws.Range("F15").Value = 135
$ws.Range("F16").Value = 136
$ws.Range("L16").Value = 'stimuli/img_uwv6y.png'
$ws.Range("M16").Value = 78.88888888888889
$ws.Range("N16").Value = 59.30555555555556
$ws.Range("O16").Value = 69.09722222222223
$ws.Range("P16").Value = 36
$ws.Range("Q16").Value = 8
$ws.Range("R16").Value = 8
$ws.Range("S16").Value = 8
$ws.Range("T16").Value = 8
$ws.Range("U16").Value = 8
$ws.Range("V16").Value = 8
$ws.Range("F17").Value = 137
$ws.Range("L17").Value = 'stimuli/img_uy1n4.png'
$ws.Range("M17").Value = 76.30555555555556
$ws.Range("N17").Value = 55.33333333333334
$ws.Range("O17").Value = 65.81944444444444
$ws.Range("P17").Value = 36
$ws.Range("Q17").Value = 7
$ws.Range("R17").Value = 7
$ws.Range("S17").Value = 7
$ws.Range("T17").Value = 7
$ws.Range("U17").Value = 7
$ws.Range("V17").Value = 7
$ws.Range("F18").Value = 138
$ws.Range("F19").Value = 139
$ws.Range("F20").Value = 140
$ws.Range("F21").Value = 141
$ws.Range("L21").Value = 'stimuli/img_jz3kd.png'
$ws.Range("M21").Value = 72.79411764705883
$ws.Range("N21").Value = 51.64705882352941
$ws.Range("O21").Value = 62.22058823529412
$ws.Range("P21").Value = 34
$ws.Range("Q21").Value = 6
$ws.Range("R21").Value = 6
$ws.Range("S21").Value = 6
$ws.Range("T21").Value = 6
$ws.Range("U21").Value = 6
$ws.Range("V21").Value = 6
$ws.Range("F22").Value = 142
$ws.Range("F23").Value = 143
$ws.Range("F24").Value = 144
$ws.Range("F25").Value = 145
$ws.Range("F26").Value = 146
$ws.Range("L26").Value = 'stimuli/img_bwo9g.png'
$ws.Range("M26").Value = 64.81818181818181
$ws.Range("N26").Value = 42.36363636363637
$ws.Range("O26").Value = 53.59090909090909
$ws.Range("P26").Value = 33
$ws.Range("Q26").Value = 4
$ws.Range("R26").Value = 4
$ws.Range("S26").Value = 4
$ws.Range("T26").Value = 4
$ws.Range("U26").Value = 4
$ws.Range("V26").Value = 4
$ws.Range("F27").Value = 147
$ws.Range("F28").Value = 148
$ws.Range("F29").Value = 149
$ws.Range("F30").Value = 150
$ws.Range("L30").Value = 'stimuli/img_ncr40.png'
$ws.Range("M30").Value = 75.66666666666667
$ws.Range("N30").Value = 54.27272727272727
$ws.Range("O30").Value = 64.96969696969697
$ws.Range("P30").Value = 33
$ws.Range("F31").Value = 151
$ws.Range("F32").Value = 152
$ws.Range("L32").Value = 'stimuli/img_c79r7.png'
$ws.Range("M32").Value = 56.26470588235294
$ws.Range("N32").Value = 34.26470588235294
$ws.Range("O32").Value = 45.26470588235294
$ws.Range("P32").Value = 34
$ws.Range("Q32").Value = 2
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 2
$ws.Range("T32").Value = 2
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 2
$ws.Range("F33").Value = 153
$ws.Range("F34").Value = 154
$ws.Range("F35").Value = 155
$ws.Range("F36").Value = 156
$ws.Range("L36").Value = 'stimuli/img_05flq.png'
$ws.Range("M36").Value = 47.10344827586207
$ws.Range("N36").Value = 25.72413793103448
$ws.Range("O36").Value = 36.41379310344828
$ws.Range("P36").Value = 29
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = 1
$ws.Range("S36").Value = 1
$ws.Range("T36").Value = 1
$ws.Range("U36").Value = 1
$ws.Range("V36").Value = 1
$ws.Range("F37").Value = 157
$ws.Range("L37").Value = 'stimuli/img_j5rpx.png'
$ws.Range("M37").Value = 72.24242424242425
$ws.Range("N37").Value = 50
$ws.Range("O37").Value = 61.12121212121212
$ws.Range("P37").Value = 33
$ws.Range("Q37").Value = 5
$ws.Range("R37").Value = 5
$ws.Range("S37").Value = 5
$ws.Range("T37").Value = 5
$ws.Range("U37").Value = 5
$ws.Range("V37").Value = 5
$ws.Range("F38").Value = 158
$ws.Range("L38").Value = 'stimuli/img_7pgd2.png'
$ws.Range("M38").Value = 78.59375
$ws.Range("N38").Value = 57.84375
$ws.Range("O38").Value = 68.21875
$ws.Range("P38").Value = 32
$ws.Range("Q38").Value = 8
$ws.Range("R38").Value = 7
$ws.Range("S38").Value = 7
$ws.Range("T38").Value = 7
$ws.Range("U38").Value = 7
$ws.Range("V38").Value = 7
$ws.Range("F39").Value = 159
$ws.Range("L39").Value = 'stimuli/img_xti0z.png'
$ws.Range("M39").Value = 81.40625
$ws.Range("N39").Value = 61.4375
$ws.Range("O39").Value = 71.421875
$ws.Range("R39").Value = 8
$ws.Range("S39").Value = 8
$ws.Range("T39").Value = 8
$ws.Range("U39").Value = 8
$ws.Range("V39").Value = 8
$ws.Range("F40").Value = 160
$ws.Range("F41").Value = 161
$ws.Range("L41").Value = 'stimuli/img_s9are.png'
$ws.Range("M41").Value = 90.14285714285714
$ws.Range("N41").Value = 75.22857142857143
$ws.Range("O41").Value = 82.68571428571428
$ws.Range("P41").Value = 35
$ws.Range("Q41").Value = 10
$ws.Range("R41").Value = 10
$ws.Range("S41").Value = 10
$ws.Range("T41").Value = 10
$ws.Range("U41").Value = 10
$ws.Range("V41").Value = 10
$ws.Range("F42").Value = 162
$ws.Range("L42").Value = 'stimuli/img_411xa.png'
$ws.Range("M42").Value = 51.03030303030303
$ws.Range("N42").Value = 28.93939393939394
$ws.Range("O42").Value = 39.98484848484848
$ws.Range("P42").Value = 33
$ws.Range("Q42").Value = 2
$ws.Range("R42").Value = 2
$ws.Range("S42").Value = 2
$ws.Range("T42").Value = 2
$ws.Range("U42").Value = 2
$ws.Range("V42").Value = 2
